# Updated cryptos list -- applies price/volume/name/link edits per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($CellRef, $Text)
    $cell = $ws.Range($CellRef)
    # Force text storage so numeric-looking strings (e.g. "560.47", "5.00")
    # are not coerced into floating point numbers, and so leading/trailing
    # spaces and non-significant zeros survive exactly as authored.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-CellText "D2" "68.963.28"
Set-CellText "E2" "  +0.66%  "
Set-CellText "D3" "2.472.90"
Set-CellText "E3" "  +0.85%  "
Set-CellText "E4" "  -0.03%  "
Set-CellText "D5" "560.47"
Set-CellText "E5" "  -0.76%  "
Set-CellText "D6" "162.26"
Set-CellText "E6" "  -0.87%  "
Set-CellText "E7" "  +0.01%  "
Set-CellText "E8" "  +0.09%  "
Set-CellText "E9" "  +0.02%  "
Set-CellText "E10" "  +0.56%  "
Set-CellText "E11" "  -2.68%  "
Set-CellText "E12" "  +1.37%  "
Set-CellText "E13" "  +0.11%  "
Set-CellText "D14" "68.862.37"
Set-CellText "E14" "  +0.66%  "
Set-CellText "E15" "  -1.27%  "
Set-CellText "D16" "23.65"
Set-CellText "E16" "  +0.15%  "
Set-CellText "D17" "2.478.35"
Set-CellText "E17" "  +0.24%  "
Set-CellText "D18" "10.73"
Set-CellText "E18" "  -2.43%  "
Set-CellText "D19" "336.26"
Set-CellText "E19" "  -2.64%  "
Set-CellText "E20" "  -2.86%  "
Set-CellText "D21" "3.79"
Set-CellText "E21" "  -0.93%  "
Set-CellText "E22" "  +0.13%  "
Set-CellText "E23" "  +0.06%  "
Set-CellText "D24" "66.76"
Set-CellText "E24" "  -2.10%  "
Set-CellText "E25" "  -2.45%  "
Set-CellText "D26" "8.21"
Set-CellText "E26" "  -0.38%  "
Set-CellText "D27" "0.0₃0818"
Set-CellText "E27" "  -2.68%  "
Set-CellText "D28" "7.22"
Set-CellText "E28" "  -1.11%  "
Set-CellText "B29" "FirstDigitalUSD"
Set-CellText "C29" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-CellText "D29" "1.00"
Set-CellText "E29" "  -0.01%  "
Set-CellText "B30" "Bittensor"
Set-CellText "C30" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-CellText "D30" "432.24"
Set-CellText "E30" "  -0.98%  "
Set-CellText "E31" "  -3.51%  "
Set-CellText "E32" "  -4.14%  "
Set-CellText "D33" "159.01"
Set-CellText "E33" "  +1.18%  "
Set-CellText "D34" "19.04"
Set-CellText "E34" "  +0.15%  "
Set-CellText "E35" "  +0.32%  "
Set-CellText "E36" "  -0.08%  "
Set-CellText "D37" "17.79"
Set-CellText "E37" "  -0.60%  "
Set-CellText "E39" "  -1.53%  "
Set-CellText "E40" "  -3.87%  "
Set-CellText "D41" "1.08"
Set-CellText "E41" "  -3.59%  "
Set-CellText "E42" "  -0.71%  "
Set-CellText "E43" "  -0.51%  "
Set-CellText "D44" "131.09"
Set-CellText "E44" "  -3.14%  "
Set-CellText "E45" "  -0.75%  "
Set-CellText "E46" "  -0.70%  "
Set-CellText "D47" "0.563"
Set-CellText "D48" "0.0913"
Set-CellText "E48" "  -0.34%  "
Set-CellText "E49" "  +0.23%  "
Set-CellText "D50" "1.39"
Set-CellText "E50" "  -2.32%  "
Set-CellText "D51" "5.00"
